# Fix start/end times: the schedule was entered 4 hours too late.
# Shift every "Starting time" (column A) and "End time (optional)" (column B)
# value in the data rows (2-32) back by 4 hours, leaving all other data
# (tables, phases, groups, team names, results) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($c in 1, 2) {
        $cell = $ws.Cells.Item($r, $c)
        $val = [string]$cell.Value2
        if ($val -match '^(\d{1,2}):(\d{2})$') {
            $h = [int]$matches[1]
            $m = $matches[2]
            $h = ($h - 4 + 24) % 24
            $cell.Value2 = ('{0:D2}:{1}' -f $h, $m)
        }
    }
}
